$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells: Wins, Losses, Ties (columns AD, AE, AF on row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the existing header row (bold, centered, top-aligned, thin border)
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$headerRange.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop
$headerRange.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$headerRange.Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# Fill in the team record (Wins/Losses/Ties) for every player row (2-51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 78   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 84   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
